$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2070063694267516
$ws.Cells.Item(2, 3).Value = 0.5445859872611465
$ws.Cells.Item(2, 10).Value = 0.01592356687898089
$ws.Cells.Item(2, 16).Value = 0.124203821656051
$ws.Cells.Item(2, 19).Value = 0.1082802547770701
$ws.Cells.Item(3, 3).Value = 0.02824858757062147
$ws.Cells.Item(3, 10).Value = 0.04519774011299435
$ws.Cells.Item(3, 16).Value = 0.7457627118644068
$ws.Cells.Item(3, 19).Value = 0.1807909604519774
$ws.Cells.Item(4, 10).Value = 0.06666666666666667
$ws.Cells.Item(4, 16).Value = 0.6222222222222222
$ws.Cells.Item(4, 19).Value = 0.3111111111111111
$ws.Cells.Item(6, 2).Value = 0.03896103896103896
$ws.Cells.Item(6, 4).Value = 0.01298701298701299
$ws.Cells.Item(6, 6).Value = 0.06926406926406926
$ws.Cells.Item(6, 10).Value = 0.303030303030303
$ws.Cells.Item(6, 15).Value = 0.01298701298701299
$ws.Cells.Item(6, 17).Value = 0.09956709956709957
$ws.Cells.Item(6, 18).Value = 0.08658008658008658
$ws.Cells.Item(6, 19).Value = 0.3766233766233766
$ws.Cells.Item(7, 2).Value = 0.108695652173913
$ws.Cells.Item(7, 4).Value = 0.005434782608695652
$ws.Cells.Item(7, 5).Value = 0.005434782608695652
$ws.Cells.Item(7, 6).Value = 0.05978260869565218
$ws.Cells.Item(7, 10).Value = 0.1413043478260869
$ws.Cells.Item(7, 15).Value = 0.02717391304347826
$ws.Cells.Item(7, 17).Value = 0.1521739130434783
$ws.Cells.Item(7, 18).Value = 0.07065217391304347
$ws.Cells.Item(7, 19).Value = 0.4293478260869565
$ws.Cells.Item(8, 2).Value = 0.1037037037037037
$ws.Cells.Item(8, 4).Value = 0.009876543209876543
$ws.Cells.Item(8, 6).Value = 0.04197530864197531
$ws.Cells.Item(8, 10).Value = 0.1407407407407407
$ws.Cells.Item(8, 15).Value = 0.007407407407407408
$ws.Cells.Item(8, 17).Value = 0.1802469135802469
$ws.Cells.Item(8, 18).Value = 0.108641975308642
$ws.Cells.Item(8, 19).Value = 0.4074074074074074
$ws.Cells.Item(9, 2).Value = 0.06796116504854369
$ws.Cells.Item(9, 4).Value = 0.01941747572815534
$ws.Cells.Item(9, 5).Value = 0.004854368932038835
$ws.Cells.Item(9, 6).Value = 0.06796116504854369
$ws.Cells.Item(9, 10).Value = 0.0970873786407767
$ws.Cells.Item(9, 15).Value = 0.02427184466019417
$ws.Cells.Item(9, 17).Value = 0.1601941747572816
$ws.Cells.Item(9, 18).Value = 0.1067961165048544
$ws.Cells.Item(9, 19).Value = 0.4514563106796117
$ws.Cells.Item(10, 2).Value = 0.1216012084592145
$ws.Cells.Item(10, 4).Value = 0.02492447129909366
$ws.Cells.Item(10, 6).Value = 0.08081570996978851
$ws.Cells.Item(10, 10).Value = 0.1057401812688822
$ws.Cells.Item(10, 15).Value = 0.01359516616314199
$ws.Cells.Item(10, 17).Value = 0.202416918429003
$ws.Cells.Item(10, 18).Value = 0.06797583081570997
$ws.Cells.Item(10, 19).Value = 0.3829305135951662
$ws.Cells.Item(11, 7).Value = 0.1153846153846154
$ws.Cells.Item(11, 10).Value = 0.09935897435897435
$ws.Cells.Item(11, 11).Value = 0.1826923076923077
$ws.Cells.Item(11, 12).Value = 0.5961538461538461
$ws.Cells.Item(11, 19).Value = 0.00641025641025641
$ws.Cells.Item(12, 7).Value = 0.6666666666666666
$ws.Cells.Item(12, 10).Value = 0.25
$ws.Cells.Item(12, 11).Value = 0.03125
$ws.Cells.Item(12, 12).Value = 0.03125
$ws.Cells.Item(12, 19).Value = 0.02083333333333333
$ws.Cells.Item(13, 7).Value = 0.6
$ws.Cells.Item(13, 10).Value = 0.32
$ws.Cells.Item(13, 19).Value = 0.08
$ws.Cells.Item(15, 6).Value = 0.01970443349753695
$ws.Cells.Item(15, 8).Value = 0.1428571428571428
$ws.Cells.Item(15, 9).Value = 0.1330049261083744
$ws.Cells.Item(15, 10).Value = 0.2906403940886699
$ws.Cells.Item(15, 11).Value = 0.07389162561576355
$ws.Cells.Item(15, 13).Value = 0.01477832512315271
$ws.Cells.Item(15, 15).Value = 0.06403940886699508
$ws.Cells.Item(15, 19).Value = 0.2610837438423645
$ws.Cells.Item(16, 6).Value = 0.02051282051282051
$ws.Cells.Item(16, 8).Value = 0.1948717948717949
$ws.Cells.Item(16, 10).Value = 0.3641025641025641
$ws.Cells.Item(16, 11).Value = 0.1128205128205128
$ws.Cells.Item(16, 13).Value = 0.03076923076923077
$ws.Cells.Item(16, 15).Value = 0.03589743589743589
$ws.Cells.Item(16, 19).Value = 0.1743589743589744
$ws.Cells.Item(17, 6).Value = 0.01658767772511848
$ws.Cells.Item(17, 8).Value = 0.1635071090047393
$ws.Cells.Item(17, 9).Value = 0.07582938388625593
$ws.Cells.Item(17, 10).Value = 0.4075829383886256
$ws.Cells.Item(17, 11).Value = 0.1042654028436019
$ws.Cells.Item(17, 13).Value = 0.02132701421800948
$ws.Cells.Item(17, 14).Value = 0.002369668246445498
$ws.Cells.Item(17, 15).Value = 0.05450236966824645
$ws.Cells.Item(17, 19).Value = 0.1540284360189574
$ws.Cells.Item(18, 6).Value = 0.02645502645502645
$ws.Cells.Item(18, 8).Value = 0.164021164021164
$ws.Cells.Item(18, 9).Value = 0.08465608465608465
$ws.Cells.Item(18, 10).Value = 0.4497354497354497
$ws.Cells.Item(18, 11).Value = 0.06349206349206349
$ws.Cells.Item(18, 13).Value = 0.01058201058201058
$ws.Cells.Item(18, 15).Value = 0.03174603174603174
$ws.Cells.Item(18, 19).Value = 0.1693121693121693
$ws.Cells.Item(19, 6).Value = 0.02001429592566119
$ws.Cells.Item(19, 8).Value = 0.174410293066476
$ws.Cells.Item(19, 9).Value = 0.08434596140100072
$ws.Cells.Item(19, 10).Value = 0.3824160114367405
$ws.Cells.Item(19, 11).Value = 0.1050750536097212
$ws.Cells.Item(19, 13).Value = 0.02072909220872051
$ws.Cells.Item(19, 14).Value = 0.0007147962830593281
$ws.Cells.Item(19, 15).Value = 0.0636168691922802
$ws.Cells.Item(19, 19).Value = 0.1486776268763402
